$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (row 1), matching the style of the
# existing header cells (bold / centered / bordered => style index 1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the look of the other header cells (bold, centered, bordered).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill season-record values for every data row (2-41): 81 wins, 81
# losses, 0 ties.
for ($row = 2; $row -le 41; $row++) {
    $ws.Cells.Item($row, 30).Value = 81
    $ws.Cells.Item($row, 31).Value = 81
    $ws.Cells.Item($row, 32).Value = 0
}
